$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$c1 = $ws.Range("A8").Characters(21, 2)
$c1.Text = "48"
$c1.Font.Name = "Andale WT"
$c1.Font.Size = 10

$c2 = $ws.Range("C9").Characters(27, 10)
$c2.Text = "11/28/2022"
$c2.Font.Name = "Andale WT"
$c2.Font.Size = 10

$c3 = $ws.Range("C9").Characters(48, 10)
$c3.Text = "12/4/2022"
$c3.Font.Name = "Andale WT"
$c3.Font.Size = 10

# --- Type-switch cells: numeric -> text "0" (copy style+value from D14) ---
$ws.Range("D14").Copy($ws.Range("C14"))
$ws.Range("D14").Copy($ws.Range("C15"))
$ws.Range("D14").Copy($ws.Range("C22"))
$ws.Range("D14").Copy($ws.Range("G22"))
$ws.Range("D14").Copy($ws.Range("C23"))
$ws.Range("D14").Copy($ws.Range("C26"))

# --- Type-switch cells: numeric -> text "***.*" (copy style+value from E15) ---
$ws.Range("E15").Copy($ws.Range("H22"))

# --- Type-switch cells: text "0" -> numeric (copy style from F26 then set value) ---
$ws.Range("F26").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("F26").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("F26").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 2
$ws.Range("F26").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("F26").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1

# --- Type-switch cells: text "***.*" -> numeric (copy style from H26 then set value) ---
$ws.Range("H26").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100
$ws.Range("H26").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -50
$ws.Range("H26").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100
$ws.Range("H26").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100

# --- Plain numeric value updates ---
$ws.Range("N14").Value = -66.666666666666
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("M15").Value = -31.578947368421
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 125
$ws.Range("F16").Value = 25
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 78.571428571428
$ws.Range("I16").Value = 204
$ws.Range("J16").Value = 161
$ws.Range("K16").Value = 26.708074534161
$ws.Range("L16").Value = 72.881355932203
$ws.Range("M16").Value = -8.108108108108
$ws.Range("N16").Value = -52.33644859813
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 23
$ws.Range("H17").Value = 21.052631578947
$ws.Range("I17").Value = 244
$ws.Range("J17").Value = 216
$ws.Range("K17").Value = 12.962962962963
$ws.Range("L17").Value = 11.415525114155
$ws.Range("M17").Value = 42.690058479532
$ws.Range("N17").Value = 10.90909090909
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -9.090909090909
$ws.Range("I18").Value = 99
$ws.Range("J18").Value = 108
$ws.Range("K18").Value = -8.333333333333
$ws.Range("L18").Value = 5.31914893617
$ws.Range("M18").Value = -65.263157894736
$ws.Range("N18").Value = -86.95652173913
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = 23.529411764705
$ws.Range("F19").Value = 67
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = 19.642857142857
$ws.Range("I19").Value = 550
$ws.Range("J19").Value = 457
$ws.Range("K19").Value = 20.35010940919
$ws.Range("L19").Value = 39.949109414758
$ws.Range("M19").Value = 21.412803532008
$ws.Range("N19").Value = 38.888888888888
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 175
$ws.Range("F20").Value = 25
$ws.Range("H20").Value = -24.242424242424
$ws.Range("I20").Value = 283
$ws.Range("J20").Value = 292
$ws.Range("K20").Value = -3.082191780821
$ws.Range("L20").Value = 99.295774647887
$ws.Range("M20").Value = 44.38775510204
$ws.Range("N20").Value = -86.215294690696
$ws.Range("C21").Value = 48
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 152
$ws.Range("G21").Value = 134
$ws.Range("H21").Value = 13.432835820895
$ws.Range("I21").Value = 1398
$ws.Range("J21").Value = 1252
$ws.Range("K21").Value = 11.661341853035
$ws.Range("L21").Value = 42.217700915564
$ws.Range("M21").Value = 3.863298662704
$ws.Range("N21").Value = -64.070932922128
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 0
$ws.Range("M23").Value = -39.215686274509
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 54.166666666666
$ws.Range("F24").Value = 103
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = -1.904761904761
$ws.Range("I24").Value = 1088
$ws.Range("J24").Value = 935
$ws.Range("K24").Value = 16.363636363636
$ws.Range("L24").Value = 28.757396449704
$ws.Range("M24").Value = -19.227913882702
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 60
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = 11.538461538461
$ws.Range("I25").Value = 424
$ws.Range("J25").Value = 403
$ws.Range("K25").Value = 5.210918114143
$ws.Range("L25").Value = 17.777777777777
$ws.Range("M25").Value = 4.176904176904
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = 50
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 47
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 17.5
$ws.Range("L27").Value = 113.636363636364
$ws.Range("G28").Value = 4
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = -13.333333333333
$ws.Range("N28").Value = -38.095238095238
$ws.Range("G29").Value = 3
$ws.Range("J29").Value = 12
$ws.Range("K29").Value = -25
$ws.Range("N29").Value = -57.142857142857
